$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column D width (closest representable value to the authored 16.6640625)
$ws.Columns.Item(4).ColumnWidth = 15.83

# Header for new column D (same bold/underline header style as A1:C1)
$ws.Range("D1").Value = "On the Schedule"
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(1, 4).PasteSpecial(-4122) | Out-Null

# New presenter rows (10-26): Last Name, First Name, youtube link, On the Schedule
$data = @(
    @("Gawronska",        "Aleksandra",    "https://youtu.be/cdqip9iDzEw"),
    @("Kloos",            "Jacob",         "https://youtu.be/fgNhWlTp3iI"),
    @("Hsing-Ming Chang", "Jamie",         "https://youtu.be/eQacyfradMY"),
    @("Galinkin",         "Ryan",          "https://www.youtube.com/watch?v=iwt4gbiIrRQ"),
    @("Bourget",          "Antione",       "https://www.youtube.com/watch?v=_mVrVeBKYTo&ab_channel=AntoineBourget"),
    @("Frantzis",         "Constantinos",  "https://youtu.be/w5VG2EgQodk"),
    @("Lolachi",          "Ramin",         "https://youtu.be/UuUktXLQMDo"),
    @("Hendrix",          "Donald",        "https://youtu.be/qcb_nn9RBgo"),
    @("Williams",         "Edward",        "https://youtu.be/xwM3ZIQfqSU"),
    @("Doner",            "Alex",          "https://youtu.be/JH0FZsQb0C8"),
    @("Ayari",            "Ethan",         "https://www.youtube.com/watch?v=e5Vtnj_TiR8"),
    @("Lino",             "Gustavo",       "https://youtu.be/juteiMQfDYg"),
    @("Halim",            "Samuel",        "https://youtu.be/dcUzcP2EDQ4"),
    @("Tolometti",        "Gavin",         "https://youtu.be/K3QhfR09egQ"),
    @("Fontes",           "Douglas",       "https://youtu.be/yrDnJWszz9g"),
    @("Trinh",            "Kevin",         "https://youtu.be/CgmLkLtnVMc"),
    @("Shackelford",      "Autum",         "https://youtu.be/yyDC75aNqVE")
)

$row = 10
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = "Y"
    $row++
}

# Also mark D2:D9 ("On the Schedule" = Y) for the pre-existing rows
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = "Y"
}

# Hyperlink for the last new row (row 26), mirroring the existing C2 hyperlink
$lastRow = 26
$ws.Hyperlinks.Add($ws.Cells.Item($lastRow, 3), "https://youtu.be/yyDC75aNqVE") | Out-Null
$ws.Cells.Item($lastRow, 3).Value = "https://youtu.be/yyDC75aNqVE"
# Re-apply the same visual style used by the pre-existing hyperlink cell (C2)
# so the new hyperlink cell matches it exactly instead of keeping the
# auto-generated style that Hyperlinks.Add applies.
$ws.Cells.Item($lastRow, 3).Style = $ws.Cells.Item(2, 3).Style

# Sheet view tweaks: select D27 (just past the new data) like the saved file
$ws.Range("D27").Select() | Out-Null
